$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 00:55"

# Row 9: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B9").Value = 13847
$ws.Range("C9").Value = 4588
$ws.Range("E9").Value = 13530
$ws.Range("G9").Value = 59
$ws.Range("H9").Value = 209

# Row 24: 'Australia' -> 'Australia'
$ws.Range("E24").Value = 703
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 7

# Row 28: 'Brasil' -> 'Brasil'
$ws.Range("B28").Value = 640
$ws.Range("C28").Value = 111
$ws.Range("E28").Value = 631

# Row 60: 'Eslovaquia' -> 'Panama'
$ws.Range("A60").Value = "Panama"
$ws.Range("B60").Value = 137
$ws.Range("C60").Value = 28
$ws.Range("D60").Value = 1
$ws.Range("E60").Value = 135
$ws.Range("F60").Value = 7
$ws.Range("H60").Value = 1

# Row 61: 'Armenia' -> 'Argentina'
$ws.Range("A61").Value = "Argentina"
$ws.Range("B61").Value = 128
$ws.Range("C61").Value = 31
$ws.Range("D61").Value = 3
$ws.Range("E61").Value = 122
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 3

# Row 62: 'Mexico' -> 'Eslovaquia'
$ws.Range("A62").Value = "Eslovaquia"
$ws.Range("B62").Value = 124
$ws.Range("C62").Value = 19
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 124
$ws.Range("F62").Value = 2
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 0

# Row 63: 'Argentina' -> 'Armenia'
$ws.Range("A63").Value = "Armenia"
$ws.Range("B63").Value = 122
$ws.Range("C63").Value = 12
$ws.Range("D63").Value = 1
$ws.Range("E63").Value = 121
$ws.Range("F63").Value = 2
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0

# Row 64: 'Croacia' -> 'Mexico'
$ws.Range("A64").Value = "Mexico"
$ws.Range("B64").Value = 118
$ws.Range("C64").Value = 25
$ws.Range("D64").Value = 4
$ws.Range("E64").Value = 113
$ws.Range("F64").Value = 1

# Row 65: 'Panama' -> 'Croacia'
$ws.Range("A65").Value = "Croacia"
$ws.Range("B65").Value = 110
$ws.Range("C65").Value = 21
$ws.Range("D65").Value = 5
$ws.Range("E65").Value = 104
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 1

# Row 118: 'Honduras' -> 'Paraguay'
$ws.Range("A118").Value = "Paraguay"
$ws.Range("B118").Value = 13
$ws.Range("C118").Value = 2
$ws.Range("E118").Value = 13
$ws.Range("F118").Value = 1

# Row 119: 'Guam' -> 'Honduras'
$ws.Range("A119").Value = "Honduras"
$ws.Range("C119").Value = 3

# Row 120: 'Nigeria' -> 'Guam'
$ws.Range("A120").Value = "Guam"
$ws.Range("D120").Value = 0
$ws.Range("E120").Value = 12

# Row 121: 'Ghana' -> 'Nigeria'
$ws.Range("A121").Value = "Nigeria"
$ws.Range("B121").Value = 12
$ws.Range("D121").Value = 1

# Row 123: 'Paraguay' -> 'Ghana'
$ws.Range("A123").Value = "Ghana"
$ws.Range("C123").Value = 4
$ws.Range("F123").Value = 0

# Row 130: 'Etiopia' -> 'Mauricio'
$ws.Range("A130").Value = "Mauricio"
$ws.Range("C130").Value = 4

# Row 131: 'Kenia' -> 'Etiopia'
$ws.Range("A131").Value = "Etiopia"
$ws.Range("C131").Value = 1

# Row 132: 'Mauricio' -> 'Kenia'
$ws.Range("A132").Value = "Kenia"
$ws.Range("C132").Value = 0

# Row 133: 'Puerto Rico' -> 'Seychelles'
$ws.Range("A133").Value = "Seychelles"

# Row 135: 'Seychelles' -> 'Puerto Rico'
$ws.Range("A135").Value = "Puerto Rico"

# Row 137: 'Tanzania' -> 'Mongolia'
$ws.Range("A137").Value = "Mongolia"
$ws.Range("C137").Value = 0

# Row 138: 'Mongolia' -> 'Tanzania'
$ws.Range("A138").Value = "Tanzania"
$ws.Range("C138").Value = 3

# Row 140: 'Aruba' -> 'Guyana'
$ws.Range("A140").Value = "Guyana"
$ws.Range("D140").Value = 0
$ws.Range("H140").Value = 1

# Row 141: 'Guyana' -> 'Aruba'
$ws.Range("A141").Value = "Aruba"
$ws.Range("D141").Value = 1
$ws.Range("H141").Value = 0
